$wb = $excel.ActiveWorkbook

# Rename sheets (task order identifiers regenerated)
$wb.Worksheets.Item("GNG_TO-16498730931363091").Name = "GNG_TO-16502912196766684"
$wb.Worksheets.Item("NB_TO-16498730945128767").Name = "NB_TO-16502912240629163"
$wb.Worksheets.Item("RS_TO-16498730945188763").Name = "RS_TO-16502912240649302"
$wb.Worksheets.Item("TOL_TO-16498730945778763").Name = "TOL_TO-16502912241419215"
$wb.Worksheets.Item("vSAT_TO-16498730946569088").Name = "vSAT_TO-16502912242369218"

# GNG sheet
$ws = $wb.Worksheets.Item("GNG_TO-16502912196766684")
$ws.Range("B2").Value = "go_stims-16502912196356695.csv"
$ws.Range("B3").Value = "GNG_stims-16502912196586702.csv"
$ws.Range("B4").Value = "go_stims-16502912196607845.csv"
$ws.Range("B5").Value = "GNG_stims-16502912196746712.csv"

# NB sheet
$ws = $wb.Worksheets.Item("NB_TO-16502912240629163")
$ws.Range("B2").Value = "OB-16502912214215589.csv"
$ws.Range("B3").Value = "OB-1650291222011692.csv"
$ws.Range("B4").Value = "ZB-match_7-16502912200905597.csv"
$ws.Range("B5").Value = "TB-16502912231204398.csv"
$ws.Range("B6").Value = "ZB-match_0-16502912199255717.csv"
$ws.Range("B7").Value = "OB-16502912215319014.csv"
$ws.Range("B8").Value = "TB-16502912223928804.csv"
$ws.Range("B9").Value = "ZB-match_2-16502912199625661.csv"
$ws.Range("B10").Value = "TB-16502912240409203.csv"

# RS sheet
$ws = $wb.Worksheets.Item("RS_TO-16502912240649302")
$ws.Range("B2").Value = "eyes open"
$ws.Range("B3").Value = "eyes closed"

# TOL sheet
$ws = $wb.Worksheets.Item("TOL_TO-16502912241419215")
$ws.Range("B2").Value = "MM_stims-16502912240929244.csv"
$ws.Range("B3").Value = "ZM_stims-16502912240679226.csv"
$ws.Range("B4").Value = "MM_stims-16502912241249225.csv"
$ws.Range("B5").Value = "ZM_stims-1650291224093934.csv"
$ws.Range("B6").Value = "MM_stims-1650291224140919.csv"
$ws.Range("B7").Value = "ZM_stims-16502912241269255.csv"

# vSAT sheet
$ws = $wb.Worksheets.Item("vSAT_TO-16502912242369218")
$ws.Range("B2").Value = "vSAT_stims-16502912242219539.csv"
$ws.Range("B3").Value = "vSAT_stims-1650291224188938.csv"
$ws.Range("B4").Value = "SAT_stims-16502912241729183.csv"
$ws.Range("B5").Value = "SAT_stims-1650291224146923.csv"
